$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.370798257884184
$ws.Range("C2").Value = -0.9727079006722509
$ws.Range("D2").Value = 0.3038797348144729
$ws.Range("E2").Value = -0.9595538525665909

$ws.Range("B3").Value = -0.08261641856586049
$ws.Range("C3").Value = -0.3163358775390812
$ws.Range("D3").Value = -0.1806755391691144
$ws.Range("E3").Value = -0.1934911213520627

$ws.Range("B4").Value = -0.8290690241082486
$ws.Range("C4").Value = -0.09790517462118323
$ws.Range("D4").Value = 0.435524577543085
$ws.Range("E4").Value = 1.013510180591872

$ws.Range("B5").Value = 0.6106966230070665
$ws.Range("C5").Value = -0.4540586176844206
$ws.Range("D5").Value = -0.5976933272505026
$ws.Range("E5").Value = 0.8935945057248388

$ws.Range("B6").Value = 1.546965732314469
$ws.Range("C6").Value = -0.3393215870801091
$ws.Range("D6").Value = 0.1707791176747497
$ws.Range("E6").Value = -1.304386256282822

$ws.Range("B7").Value = -0.8093931452109928
$ws.Range("C7").Value = -1.090172000347444
$ws.Range("D7").Value = -0.4531575605523151
$ws.Range("E7").Value = -1.792323026094754
